$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.593000000000001

$ws.Range("B7").Value = 5.13
$ws.Range("C7").Value = -13.208

$ws.Range("C15").Value = -13.173

$ws.Range("B16").Value = 5.456999999999999
$ws.Range("D16").Value = -8.614000000000001

$ws.Range("D19").Value = -7.742

$ws.Range("C21").Value = -12.132

$ws.Range("C22").Value = -12.961

$ws.Range("C23").Value = -12.223

$ws.Range("B28").Value = 5.709000000000001

$ws.Range("B29").Value = 5.065

$ws.Range("B32").Value = 6.648000000000001

$ws.Range("C34").Value = -11.997

$ws.Range("D36").Value = -7.625

$ws.Range("B40").Value = 9.044999999999998

$ws.Range("C43").Value = -13.516

$ws.Range("C45").Value = -13.03

$ws.Range("D46").Value = -8.372

$ws.Range("C50").Value = -13.018
$ws.Range("D50").Value = -8.385999999999999

$ws.Range("C51").Value = -11.276

$ws.Range("B52").Value = 5.486

$ws.Range("B57").Value = 5.093000000000001

$ws.Range("B66").Value = 5.054
$ws.Range("C66").Value = -10.883

$ws.Range("C67").Value = -11.395

$ws.Range("C79").Value = -12.117

$ws.Range("C84").Value = -14.098

$ws.Range("C92").Value = -11.139

$ws.Range("D95").Value = -7.903999999999999

$ws.Range("C97").Value = -12.146
$ws.Range("D97").Value = -8.599

$ws.Range("B100").Value = 5.939000000000001
